# Updated state-transition probability matrix after simulating additional games.
# Row/col headers (Starting_State, Af0..Br0) are left untouched; only the
# recomputed transition probabilities (one row per starting state) are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updatedProbabilities = @{
    2 = @{
        "B" = 0.2084805653710247
        "C" = 0.5265017667844523
        "J" = 0.01413427561837456
        "P" = 0.1519434628975265
        "S" = 0.0989399293286219
    }
    3 = @{
        "B" = 0.006451612903225806
        "C" = 0.03870967741935484
        "J" = 0.05806451612903226
        "P" = 0.7483870967741936
        "S" = 0.1483870967741935
    }
    4 = @{
        "J" = 0.02380952380952381
        "P" = 0.7857142857142857
        "S" = 0.1904761904761905
    }
    5 = @{
        "J" = 0.5
        "P" = 0.5
    }
    6 = @{
        "B" = 0.06060606060606061
        "D" = 0.02164502164502164
        "E" = 0.004329004329004329
        "F" = 0.0303030303030303
        "J" = 0.4025974025974026
        "O" = 0.01298701298701299
        "Q" = 0.1385281385281385
        "R" = 0.05627705627705628
        "S" = 0.2727272727272727
    }
    7 = @{
        "B" = 0.08888888888888889
        "F" = 0.02777777777777778
        "J" = 0.1055555555555556
        "O" = 0.01666666666666667
        "Q" = 0.2333333333333333
        "R" = 0.08888888888888889
        "S" = 0.4055555555555556
    }
    8 = @{
        "B" = 0.1006289308176101
        "D" = 0.01886792452830189
        "E" = 0.00419287211740042
        "F" = 0.05660377358490566
        "J" = 0.1425576519916142
        "O" = 0.01467505241090147
        "Q" = 0.1761006289308176
        "R" = 0.0880503144654088
        "S" = 0.3983228511530398
    }
    9 = @{
        "B" = 0.08074534161490683
        "D" = 0.02484472049689441
        "F" = 0.03726708074534162
        "J" = 0.1180124223602484
        "O" = 0.03726708074534162
        "Q" = 0.1801242236024845
        "R" = 0.1055900621118012
        "S" = 0.4161490683229814
    }
    10 = @{
        "B" = 0.09553158705701079
        "D" = 0.01540832049306626
        "F" = 0.0600924499229584
        "J" = 0.1063174114021572
        "O" = 0.01771956856702619
        "Q" = 0.2326656394453005
        "R" = 0.1016949152542373
        "S" = 0.3705701078582435
    }
    11 = @{
        "G" = 0.1404682274247492
        "J" = 0.1103678929765886
        "K" = 0.1906354515050167
        "L" = 0.5418060200668896
        "S" = 0.01672240802675585
    }
    12 = @{
        "G" = 0.7289156626506024
        "J" = 0.1927710843373494
        "K" = 0.01204819277108434
        "L" = 0.006024096385542169
        "S" = 0.06024096385542169
    }
    13 = @{
        "G" = 0.6904761904761905
        "J" = 0.2142857142857143
        "S" = 0.09523809523809523
    }
    15 = @{
        "F" = 0.02369668246445497
        "H" = 0.1516587677725119
        "I" = 0.06161137440758294
        "J" = 0.3412322274881517
        "K" = 0.04739336492890995
        "M" = 0.02369668246445497
        "O" = 0.06161137440758294
        "S" = 0.2890995260663507
    }
    16 = @{
        "F" = 0.04232804232804233
        "H" = 0.1904761904761905
        "I" = 0.0582010582010582
        "J" = 0.3915343915343915
        "K" = 0.1111111111111111
        "M" = 0.01587301587301587
        "O" = 0.07407407407407407
        "S" = 0.1164021164021164
    }
    17 = @{
        "F" = 0.0389344262295082
        "H" = 0.2069672131147541
        "I" = 0.07377049180327869
        "J" = 0.3770491803278688
        "K" = 0.1065573770491803
        "M" = 0.01434426229508197
        "O" = 0.05532786885245902
        "S" = 0.1270491803278689
    }
    18 = @{
        "F" = 0.02232142857142857
        "H" = 0.2008928571428572
        "I" = 0.09375
        "J" = 0.46875
        "K" = 0.08035714285714286
        "M" = 0.008928571428571428
        "O" = 0.05803571428571429
        "S" = 0.06696428571428571
    }
    19 = @{
        "F" = 0.03092783505154639
        "H" = 0.2117367168913561
        "I" = 0.06344171292624901
        "J" = 0.359238699444885
        "K" = 0.1046788263283109
        "M" = 0.02061855670103093
        "O" = 0.06344171292624901
        "S" = 0.1459159397303727
    }
}

foreach ($row in $updatedProbabilities.Keys) {
    $rowValues = $updatedProbabilities[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
